$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.726.42'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').Value = '2.733.05'
$ws.Range('E3').Value = '  -0.56%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'562.34"
$ws.Range('E5').Value = '  -2.03%  '
$ws.Range('D6').Value = "'159.87"
$ws.Range('E6').Value = '  +1.71%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -0.50%  '
$ws.Range('E9').Value = '  +0.22%  '
$ws.Range('E10').Value = '  +4.02%  '
$ws.Range('E11').Value = '  +2.48%  '
$ws.Range('E12').Value = '  -0.64%  '
$ws.Range('D13').Value = '3.215.90'
$ws.Range('E13').Value = '  -0.50%  '
$ws.Range('E14').Value = '  +1.65%  '
$ws.Range('D15').Value = '63.548.09'
$ws.Range('E15').Value = '  -0.04%  '
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').Value = '2.735.90'
$ws.Range('D18').Value = "'12.54"
$ws.Range('E18').Value = '  +3.42%  '
$ws.Range('E19').Value = '  -0.83%  '
$ws.Range('D20').Value = "'354.43"
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('D21').Value = "'6.56"
$ws.Range('E21').Value = '  -2.72%  '
$ws.Range('E22').Value = '  +0.26%  '
$ws.Range('D23').Value = "'0.521"
$ws.Range('E23').Value = '  -2.95%  '
$ws.Range('D24').Value = "'64.38"
$ws.Range('E24').Value = '  -1.24%  '
$ws.Range('D25').Value = "'0.170"
$ws.Range('E25').Value = '  +0.41%  '
$ws.Range('D26').Value = "'0.999"
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').Value = "'8.38"
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('E28').Value = '  +1.17%  '
$ws.Range('E29').Value = '  +1.15%  '
$ws.Range('D30').Value = "'7.20"
$ws.Range('E30').Value = '  +3.72%  '
$ws.Range('E31').Value = '  +10.01%  '
$ws.Range('D32').Value = "'165.80"
$ws.Range('E32').Value = '  -2.15%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').Value = "'4.89"
$ws.Range('E33').Value = '  +1.16%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = "'20.03"
$ws.Range('E34').Value = '  -0.43%  '
$ws.Range('E36').Value = '  +2.31%  '
$ws.Range('E37').Value = '  +1.06%  '
$ws.Range('D38').Value = "'0.974"
$ws.Range('E38').Value = '  -0.60%  '
$ws.Range('D39').Value = "'346.38"
$ws.Range('E39').Value = '  +6.51%  '
$ws.Range('D40').Value = "'6.28"
$ws.Range('E40').Value = '  +1.99%  '
$ws.Range('D41').Value = "'4.10"
$ws.Range('E41').Value = '  -0.65%  '
$ws.Range('D42').Value = "'38.54"
$ws.Range('E42').Value = '  -0.90%  '
$ws.Range('D43').Value = "'21.87"
$ws.Range('E43').Value = '  +2.89%  '
$ws.Range('D44').Value = "'21.05"
$ws.Range('E44').Value = '  -0.75%  '
$ws.Range('E45').Value = '  -0.44%  '
$ws.Range('E46').Value = '  +0.86%  '
$ws.Range('E47').Value = '  -1.34%  '
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('D50').Value = "'131.80"
$ws.Range('E50').Value = '  -2.07%  '
$ws.Range('D51').Value = "'11.07"
$ws.Range('E51').Value = '  +0.13%  '
